$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Nenkulm" floor ("basement-below") to "Kati1" ("Floor1")
# throughout the second table (rows 21-40) in the "Etazhi" (floor) column F,
# and in the merged label cells G21 / H21.

$nenkulmRows = 21..34
foreach ($r in $nenkulmRows) {
    $ws.Range("F$r").Value = "Kati1"
}

$nenkulmTerasaRows = 35..40
foreach ($r in $nenkulmTerasaRows) {
    $ws.Range("F$r").Value = "Kati1-TERASË"
}

$ws.Range("G21").Value = "Kati1"
$ws.Range("H21").Value = "Kati1"

# Update the sheet selection / view to match the edited location.
$null = $ws.Range("B20:J40").Select()
